# Natmi following Dr Hou advice
# Re-run of the NATMI LR-pair computation for Rspo3-Lgr4 (FAPs sender) now that a
# new target cluster "M2" has been added alongside the existing ECs/FAPs/M1/Neutro/sCs
# clusters. This refreshes the existing data rows and appends a new row for the
# M2 target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $values[$i]
    }
}

# Row 2: Target cluster ECs
Set-RowValues 2 @("FAPs","Rspo3","Lgr4","ECs",3,1,2.822099333333334,8.466298,1,1,3,1,1.263917333333333,3.791752,0.06307655211196754,0.06307655211196754,3.566900263788445,32.102102374096,0.06307655211196754,0.06307655211196754)

# Row 3: Target cluster FAPs
Set-RowValues 3 @("FAPs","Rspo3","Lgr4","FAPs",3,1,2.822099333333334,8.466298,1,1,3,1,10.871597,32.614791,0.542553564719536,0.5425535647195359,30.68072664596867,276.126539813718,0.542553564719536,0.5425535647195359)

# Row 4: Target cluster M1
Set-RowValues 4 @("FAPs","Rspo3","Lgr4","M1",3,1,2.822099333333334,8.466298,1,1,3,1,0.1440293333333333,0.432088,0.007187870211173049,0.007187870211173049,0.4064650855804445,3.658185770224,0.007187870211173049,0.007187870211173049)

# Row 5: Target cluster M2 (newly introduced cluster)
Set-RowValues 5 @("FAPs","Rspo3","Lgr4","M2",3,1,2.822099333333334,8.466298,1,1,1,0.3333333333333333,0.114435,0.343305,0.005710947267331571,0.00571094726733157,0.3229469372100001,2.90652243489,0.005710947267331571,0.00571094726733157)

# Row 6: Target cluster Neutro
Set-RowValues 6 @("FAPs","Rspo3","Lgr4","Neutro",3,1,2.822099333333334,8.466298,1,1,3,1,0.6215056666666667,1.864517,0.03101661282545625,0.03101661282545625,1.753950727562889,15.785556548066,0.03101661282545625,0.03101661282545625)

# Row 7 (new): Target cluster sCs
Set-RowValues 7 @("FAPs","Rspo3","Lgr4","sCs",3,1,2.822099333333334,8.466298,1,1,3,1,7.022347333333333,21.067042,0.3504544528645357,0.3504544528645357,19.81776172783511,178.359855550516,0.3504544528645357,0.3504544528645357)
